$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1924.5
$ws.Range("I2").Value = 3000
$ws.Range("K2").Value = 3000
$ws.Range("M2").Value = -2887
$ws.Range("H4").Value = 237.71428
$ws.Range("I4").Value = 254
$ws.Range("J4").Value = 140
$ws.Range("K4").Value = 254
$ws.Range("L4").Value = 140
$ws.Range("M4").Value = -140
$ws.Range("N4").Value = -368
$ws.Range("H5").Value = 44.57143
$ws.Range("I5").Value = 44.57143
$ws.Range("K5").Value = 44.57143
$ws.Range("M5").Value = 70.42857000000001
$ws.Range("H6").Value = 820.35
$ws.Range("I6").Value = 616.0625
$ws.Range("K6").Value = 1848.1875
$ws.Range("M6").Value = -1736.1875
$ws.Range("H12").Value = 248.83333
$ws.Range("I12").Value = 199.6
$ws.Range("K12").Value = 199.6
$ws.Range("M12").Value = -29.59999999999999
$ws.Range("H13").Value = 9999.5
$ws.Range("J13").Value = 9999.5
$ws.Range("L13").Value = 9999.5
$ws.Range("N13").Value = -10337.5
$ws.Range("H19").Value = 3824
$ws.Range("I19").Value = 5159.2
$ws.Range("K19").Value = 5159.2
$ws.Range("M19").Value = -4984.2
$ws.Range("H29").Value = 10000.25
$ws.Range("I29").Value = 1003
$ws.Range("J29").Value = 12999.333
$ws.Range("K29").Value = 3009
$ws.Range("L29").Value = 38997.999
$ws.Range("M29").Value = -2728
$ws.Range("N29").Value = -39559.999
$ws.Range("H32").Value = 14462.1875
$ws.Range("J32").Value = 15626.5
$ws.Range("L32").Value = 15626.5
$ws.Range("N32").Value = -16278.5
$ws.Range("H33").Value = 312
$ws.Range("I33").Value = 164
$ws.Range("K33").Value = 164
$ws.Range("M33").Value = 65
$ws.Range("H38").Value = 25.25
$ws.Range("I38").Value = 25.25
$ws.Range("K38").Value = 75.75
$ws.Range("M38").Value = 296.25
$ws.Range("H39").Value = 499
$ws.Range("I39").Value = 498
$ws.Range("K39").Value = 1494
$ws.Range("M39").Value = -1198
$ws.Range("H48").Value = 6000
$ws.Range("I48").Value = 6090.909
$ws.Range("K48").Value = 18272.727
$ws.Range("M48").Value = -17980.727
$ws.Range("H51").Value = 7819412
$ws.Range("J51").Value = 8009.25
$ws.Range("L51").Value = 8009.25
$ws.Range("N51").Value = -8977.25
$ws.Range("H56").Value = 6000
$ws.Range("I56").Value = 6090.909
$ws.Range("K56").Value = 18272.727
$ws.Range("M56").Value = -17738.727
$ws.Range("H58").Value = 5080.6665
$ws.Range("I58").Value = 1042.5
$ws.Range("J58").Value = 7099.75
$ws.Range("K58").Value = 3127.5
$ws.Range("L58").Value = 21299.25
$ws.Range("M58").Value = -2977.5
$ws.Range("N58").Value = -21599.25
$ws.Range("H64").Value = 26850.1
$ws.Range("J64").Value = 9999.5
$ws.Range("L64").Value = 9999.5
$ws.Range("N64").Value = -10495.5
$ws.Range("H67").Value = 26850.1
$ws.Range("J67").Value = 9999.5
$ws.Range("L67").Value = 9999.5
$ws.Range("N67").Value = -11715.5
$ws.Range("H69").Value = 178474.8
$ws.Range("I69").Value = 23332
$ws.Range("J69").Value = 205852.94
$ws.Range("K69").Value = 69996
$ws.Range("L69").Value = 617558.8200000001
$ws.Range("M69").Value = -69122
$ws.Range("N69").Value = -619306.8200000001
$ws.Range("H72").Value = 178474.8
$ws.Range("I72").Value = 23332
$ws.Range("J72").Value = 205852.94
$ws.Range("K72").Value = 209988
$ws.Range("L72").Value = 1852676.46
$ws.Range("M72").Value = -205620
$ws.Range("N72").Value = -1861412.46
$ws.Range("H80").Value = 8496.333000000001
$ws.Range("I80").Value = 1020.4286
$ws.Range("J80").Value = 15037.75
$ws.Range("K80").Value = 3061.2858
$ws.Range("L80").Value = 45113.25
$ws.Range("M80").Value = -2063.2858
$ws.Range("N80").Value = -47109.25
$ws.Range("H83").Value = 8496.333000000001
$ws.Range("I83").Value = 1020.4286
$ws.Range("J83").Value = 15037.75
$ws.Range("K83").Value = 9183.857399999999
$ws.Range("L83").Value = 135339.75
$ws.Range("M83").Value = -4191.857399999999
$ws.Range("N83").Value = -145323.75
$ws.Range("H98").Value = 4726.647
$ws.Range("I98").Value = 2689.4666
$ws.Range("K98").Value = 2689.4666
$ws.Range("M98").Value = -1191.4666
$ws.Range("H101").Value = 552.6
$ws.Range("J101").Value = 2000
$ws.Range("L101").Value = 6000
$ws.Range("N101").Value = -9244
$ws.Range("H107").Value = 1241.2727
$ws.Range("I107").Value = 850.17645
$ws.Range("K107").Value = 850.17645
$ws.Range("M107").Value = 1069.82355
$ws.Range("H113").Value = 6706.5
$ws.Range("I113").Value = 6301.357
$ws.Range("K113").Value = 6301.357
$ws.Range("M113").Value = -3047.357
$ws.Range("H122").Value = 4726.647
$ws.Range("I122").Value = 2689.4666
$ws.Range("K122").Value = 8068.399800000001
$ws.Range("M122").Value = -5618.399800000001
$ws.Range("H132").Value = 3191556.5
$ws.Range("I132").Value = 3502810.8
$ws.Range("K132").Value = 10508432.4
$ws.Range("M132").Value = -10505902.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19411.373
$ws.Range("I32").Value = 20331.236
$ws.Range("K32").Value = 20331.236
$ws.Range("M32").Value = -20044.236
$ws.Range("H45").Value = 3545.4285
$ws.Range("I45").Value = 1415.6666
$ws.Range("J45").Value = 5142.75
$ws.Range("K45").Value = 1415.6666
$ws.Range("L45").Value = 5142.75
$ws.Range("M45").Value = -1038.6666
$ws.Range("N45").Value = -5896.75
$ws.Range("H61").Value = 52500
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H63").Value = 3326.6875
$ws.Range("I63").Value = 3053.2222
$ws.Range("J63").Value = 3678.2856
$ws.Range("K63").Value = 3053.2222
$ws.Range("L63").Value = 3678.2856
$ws.Range("M63").Value = -2367.2222
$ws.Range("N63").Value = -5050.2856
$ws.Range("H66").Value = 3326.6875
$ws.Range("I66").Value = 3053.2222
$ws.Range("J66").Value = 3678.2856
$ws.Range("K66").Value = 15266.111
$ws.Range("L66").Value = 18391.428
$ws.Range("M66").Value = -11834.111
$ws.Range("N66").Value = -25255.428
$ws.Range("H74").Value = 471564
$ws.Range("I74").Value = 1000917.3
$ws.Range("J74").Value = 17832.572
$ws.Range("K74").Value = 1000917.3
$ws.Range("L74").Value = 17832.572
$ws.Range("M74").Value = -1000043.3
$ws.Range("N74").Value = -19580.572
$ws.Range("H77").Value = 471564
$ws.Range("I77").Value = 1000917.3
$ws.Range("J77").Value = 17832.572
$ws.Range("K77").Value = 5004586.5
$ws.Range("L77").Value = 89162.86
$ws.Range("M77").Value = -5000218.5
$ws.Range("N77").Value = -97898.86
$ws.Range("H102").Value = 2723.077
$ws.Range("I102").Value = 2649.2
$ws.Range("K102").Value = 2649.2
$ws.Range("M102").Value = -1027.2
$ws.Range("H110").Value = 38019
$ws.Range("I110").Value = 42471.477
$ws.Range("K110").Value = 42471.477
$ws.Range("M110").Value = -40426.477
$ws.Range("H122").Value = 1450.931
$ws.Range("I122").Value = 1242.72
$ws.Range("K122").Value = 3728.16
$ws.Range("M122").Value = -1278.16
$ws.Range("H132").Value = 3325.6
$ws.Range("I132").Value = 2928.6667
$ws.Range("K132").Value = 8786.000100000001
$ws.Range("M132").Value = -6256.000100000001
$ws.Range("H136").Value = 52500
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2331.2424
$ws.Range("I107").Value = 2222.8386
$ws.Range("K107").Value = 2222.8386
$ws.Range("M107").Value = -302.8386
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2230.2144
$ws.Range("I16").Value = 2101.9167
$ws.Range("K16").Value = 2101.9167
$ws.Range("M16").Value = -1814.9167
$ws.Range("H31").Value = 6667604
$ws.Range("J31").Value = 899
$ws.Range("L31").Value = 899
$ws.Range("N31").Value = -1489
$ws.Range("H34").Value = 6667604
$ws.Range("J34").Value = 899
$ws.Range("L34").Value = 899
$ws.Range("N34").Value = -1303
$ws.Range("H58").Value = 1881.7142
$ws.Range("I58").Value = 1100.8125
$ws.Range("J58").Value = 4380.6
$ws.Range("K58").Value = 1100.8125
$ws.Range("L58").Value = 4380.6
$ws.Range("M58").Value = -897.8125
$ws.Range("N58").Value = -4786.6
$ws.Range("H62").Value = 4283.5
$ws.Range("I62").Value = 4098.3335
$ws.Range("K62").Value = 4098.3335
$ws.Range("M62").Value = -3474.3335
$ws.Range("H65").Value = 4283.5
$ws.Range("I65").Value = 4098.3335
$ws.Range("K65").Value = 20491.6675
$ws.Range("M65").Value = -17371.6675
$ws.Range("H105").Value = 1650.1
$ws.Range("I105").Value = 812.75
$ws.Range("J105").Value = 4999.5
$ws.Range("K105").Value = 812.75
$ws.Range("L105").Value = 4999.5
$ws.Range("M105").Value = 934.25
$ws.Range("N105").Value = -8493.5
$ws.Range("H107").Value = 700.8333
$ws.Range("I107").Value = 507.66666
$ws.Range("J107").Value = 1666.6666
$ws.Range("K107").Value = 507.66666
$ws.Range("L107").Value = 1666.6666
$ws.Range("M107").Value = 1412.33334
$ws.Range("N107").Value = -5506.6666
$ws.Range("H113").Value = 2230.2144
$ws.Range("I113").Value = 2101.9167
$ws.Range("K113").Value = 2101.9167
$ws.Range("M113").Value = 68.08329999999978
$ws.Range("H115").Value = 30356.785
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -32350
$ws.Range("H122").Value = 1568.2
$ws.Range("I122").Value = 1644.2307
$ws.Range("K122").Value = 4932.6921
$ws.Range("M122").Value = -2482.6921
$ws.Range("H132").Value = 38107.215
$ws.Range("I132").Value = 60812
$ws.Range("J132").Value = 3018
$ws.Range("K132").Value = 182436
$ws.Range("L132").Value = 9054
$ws.Range("M132").Value = -179906
$ws.Range("N132").Value = -14114
$ws.Range("H134").Value = 4128.1
$ws.Range("I134").Value = 3696.6
$ws.Range("K134").Value = 11089.8
$ws.Range("M134").Value = -8554.799999999999
$ws.Range("H136").Value = 1881.7142
$ws.Range("I136").Value = 1100.8125
$ws.Range("J136").Value = 4380.6
$ws.Range("K136").Value = 3302.4375
$ws.Range("L136").Value = 13141.8
$ws.Range("M136").Value = -752.4375
$ws.Range("N136").Value = -18241.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 262
$ws.Range("I2").Value = 227.5
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 1365
$ws.Range("L2").Value = 2400
$ws.Range("M2").Value = -1252
$ws.Range("N2").Value = -2626
$ws.Range("H56").Value = 6310.25
$ws.Range("I56").Value = 6310.25
$ws.Range("K56").Value = 6310.25
$ws.Range("M56").Value = -5780.25
$ws.Range("H136").Value = 2106.6667
$ws.Range("I136").Value = 2106.6667
$ws.Range("K136").Value = 6320.000100000001
$ws.Range("M136").Value = -1220.000100000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9266.066000000001
$ws.Range("I80").Value = 3910.111
$ws.Range("J80").Value = 17300
$ws.Range("K80").Value = 3910.111
$ws.Range("L80").Value = 17300
$ws.Range("M80").Value = -2912.111
$ws.Range("N80").Value = -19296
$ws.Range("H83").Value = 9266.066000000001
$ws.Range("I83").Value = 3910.111
$ws.Range("J83").Value = 17300
$ws.Range("K83").Value = 19550.555
$ws.Range("L83").Value = 86500
$ws.Range("M83").Value = -14558.555
$ws.Range("N83").Value = -96484
$ws.Range("H102").Value = 17922
$ws.Range("I102").Value = 20208.076
$ws.Range("J102").Value = 3062.5
$ws.Range("K102").Value = 20208.076
$ws.Range("L102").Value = 3062.5
$ws.Range("M102").Value = -18586.076
$ws.Range("N102").Value = -6306.5
$ws.Range("H122").Value = 2647.946
$ws.Range("I122").Value = 2341
$ws.Range("K122").Value = 7023
$ws.Range("M122").Value = -4573
$ws.Range("H132").Value = 2142.375
$ws.Range("I132").Value = 1947.4286
$ws.Range("J132").Value = 3507
$ws.Range("K132").Value = 5842.2858
$ws.Range("L132").Value = 10521
$ws.Range("M132").Value = -3312.2858
$ws.Range("N132").Value = -15581
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4420.222
$ws.Range("I7").Value = 5166.6665
$ws.Range("K7").Value = 5166.6665
$ws.Range("M7").Value = -5054.6665
$ws.Range("H16").Value = 583.6667
$ws.Range("J16").Value = 676
$ws.Range("L16").Value = 676
$ws.Range("N16").Value = -1016
$ws.Range("H22").Value = 1452.2916
$ws.Range("I22").Value = 1155.1818
$ws.Range("J22").Value = 1703.6923
$ws.Range("K22").Value = 1155.1818
$ws.Range("L22").Value = 1703.6923
$ws.Range("M22").Value = -860.1818000000001
$ws.Range("N22").Value = -2293.6923
$ws.Range("H27").Value = 1452.2916
$ws.Range("I27").Value = 1155.1818
$ws.Range("J27").Value = 1703.6923
$ws.Range("K27").Value = 1155.1818
$ws.Range("L27").Value = 1703.6923
$ws.Range("M27").Value = -1048.1818
$ws.Range("N27").Value = -1917.6923
$ws.Range("H40").Value = 4856.4287
$ws.Range("I40").Value = 4856.4287
$ws.Range("K40").Value = 4856.4287
$ws.Range("M40").Value = -4720.4287
$ws.Range("H55").Value = 2024.5294
$ws.Range("I55").Value = 1179.1818
$ws.Range("J55").Value = 3574.3333
$ws.Range("K55").Value = 1179.1818
$ws.Range("L55").Value = 3574.3333
$ws.Range("M55").Value = -1006.1818
$ws.Range("N55").Value = -3920.3333
$ws.Range("H82").Value = 1896.5454
$ws.Range("J82").Value = 2229.6667
$ws.Range("L82").Value = 2229.6667
$ws.Range("N82").Value = -2951.6667
$ws.Range("H85").Value = 1896.5454
$ws.Range("J85").Value = 2229.6667
$ws.Range("L85").Value = 2229.6667
$ws.Range("N85").Value = -4725.6667
$ws.Range("H93").Value = 2413.7144
$ws.Range("I93").Value = 1158.5
$ws.Range("J93").Value = 5551.75
$ws.Range("K93").Value = 1158.5
$ws.Range("L93").Value = 5551.75
$ws.Range("M93").Value = 89.5
$ws.Range("N93").Value = -8047.75
$ws.Range("H126").Value = 4420.222
$ws.Range("I126").Value = 5166.6665
$ws.Range("K126").Value = 15499.9995
$ws.Range("M126").Value = -13029.9995
$ws.Range("H132").Value = 8425.286
$ws.Range("J132").Value = 5744.75
$ws.Range("L132").Value = 17234.25
$ws.Range("N132").Value = -22294.25
$ws.Range("H134").Value = 104998.336
$ws.Range("J134").Value = 104998.336
$ws.Range("L134").Value = 104998.336
$ws.Range("N134").Value = -115138.336
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 17592.818
$ws.Range("J45").Value = 17592.818
$ws.Range("L45").Value = 17592.818
$ws.Range("N45").Value = -18574.818
$ws.Range("H81").Value = 8291
$ws.Range("I81").Value = 9665.777
$ws.Range("K81").Value = 19331.554
$ws.Range("M81").Value = -18270.554
$ws.Range("H84").Value = 8291
$ws.Range("I84").Value = 9665.777
$ws.Range("K84").Value = 96657.77
$ws.Range("M84").Value = -91353.77
$ws.Range("H93").Value = 79333.336
$ws.Range("J93").Value = 79333.336
$ws.Range("L93").Value = 79333.336
$ws.Range("N93").Value = -84325.336
$ws.Range("H107").Value = 1417.25
$ws.Range("J107").Value = 1135
$ws.Range("L107").Value = 3405
$ws.Range("N107").Value = -7245
$ws.Range("H113").Value = 1116.8125
$ws.Range("I113").Value = 1023.4545
$ws.Range("J113").Value = 1322.2
$ws.Range("K113").Value = 3070.3635
$ws.Range("L113").Value = 3966.6
$ws.Range("M113").Value = -900.3635000000004
$ws.Range("N113").Value = -8306.6
$ws.Range("H132").Value = 119566.5
$ws.Range("I132").Value = 234166.67
$ws.Range("J132").Value = 4966.3335
$ws.Range("K132").Value = 702500.01
$ws.Range("L132").Value = 14899.0005
$ws.Range("M132").Value = -699970.01
$ws.Range("N132").Value = -19959.0005
